$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 12
$ws.Cells.Item(1, 2).NumberFormat = "@"
$ws.Cells.Item(1, 2).Value = '2009-08-08'
$ws.Cells.Item(1, 2).ClearFormats()
$ws.Cells.Item(1, 3).Value = 12
$ws.Cells.Item(1, 4).Value = 'Jgbcyf'

$ws.Cells.Item(2, 1).Value = 12
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = '2008-07-09'
$ws.Cells.Item(2, 2).ClearFormats()
$ws.Cells.Item(2, 3).Value = 321
$ws.Cells.Item(2, 4).Value = 'Описание ошибки'

$ws.Cells.Item(3, 1).Value = 12
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = '2003-01-01'
$ws.Cells.Item(3, 2).ClearFormats()
$ws.Cells.Item(3, 3).Value = 123
$ws.Cells.Item(3, 4).Value = 'Ошибка энкодера поворота лопасти'

$ws.Cells.Item(4, 1).Value = 12
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = '2001-01-01'
$ws.Cells.Item(4, 2).ClearFormats()
$ws.Cells.Item(4, 3).Value = 345
$ws.Cells.Item(4, 4).Value = 'Ошибка главного контроллера'

$ws.Cells.Item(5, 1).Value = 12
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = '2001-01-01'
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 3).Value = 123
$ws.Cells.Item(5, 4).Value = 'Ошибка чегото чегото'

$ws.Cells.Item(6, 1).Value = 22
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = '2009-04-05'
$ws.Cells.Item(6, 2).ClearFormats()
$ws.Cells.Item(6, 3).Value = 123
$ws.Cells.Item(6, 4).Value = 'aljhgqalf'

$ws.Cells.Item(7, 1).Value = 22
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = '2008-09-09'
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).Value = 44
$ws.Cells.Item(7, 4).Value = 'Описание ошибки'

$ws.Cells.Item(8, 1).Value = 22
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = '2008-04-05'
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).Value = 453
$ws.Cells.Item(8, 4).Value = 'Описание ошибки рррр'

$ws.Cells.Item(9, 1).Value = 22
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = '2001-01-01'
$ws.Cells.Item(9, 2).ClearFormats()
$ws.Cells.Item(9, 3).Value = 123
$ws.Cells.Item(9, 4).Value = 'Ошибка ошибка ошибка'

$ws.Cells.Item(10, 1).Value = 2222
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = '2024-01-06'
$ws.Cells.Item(10, 2).ClearFormats()
$ws.Cells.Item(10, 3).Value = 245
$ws.Cells.Item(10, 4).Value = 'gjdgfxbg'

$ws.Cells.Item(11, 1).Value = 314
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = '2001-01-01'
$ws.Cells.Item(11, 2).ClearFormats()
$ws.Cells.Item(11, 3).Value = 34
$ws.Cells.Item(11, 4).Value = 'jgbcyfb ороплы лдоцке'

$ws.Cells.Item(12, 1).Value = 33
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = '2008-09-05'
$ws.Cells.Item(12, 2).ClearFormats()
$ws.Cells.Item(12, 3).Value = 134
$ws.Cells.Item(12, 4).Value = 'Описание ошибки'

$ws.Cells.Item(13, 1).Value = 442
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = '2001-05-07'
$ws.Cells.Item(13, 2).ClearFormats()
$ws.Cells.Item(13, 3).Value = 425
$ws.Cells.Item(13, 4).Value = 'Ошибкак главного вала выскокая темепратруа'

$ws.Cells.Item(14, 1).Value = 555
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = '2099-01-01'
$ws.Cells.Item(14, 2).ClearFormats()
$ws.Cells.Item(14, 3).Value = 345
$ws.Cells.Item(14, 4).Value = 'fsdghsh'
